# #69 Layer.copyRect, operateRect, stretchCopy, operateStretch, affineCopy,
# operateAffine: allow Bitmap class (not just Layer) for `src`. Bump version.
#
# This adds two new message rows (TVPSpecifyLayerOrBitmap /
# TVPCannotAcceptModeAuto) to the "messages" table on the "コア共通"
# worksheet, inserted just above the existing TVPCannotCreateEmptyLayerImage
# row, and moves the active sheet/selection onto that worksheet.

$wb = $excel.ActiveWorkbook

$wsCore = $wb.Worksheets.Item(2)   # "コア共通" (Core common)
$wsWin32 = $wb.Worksheets.Item(3)  # "コアWin32" (Core Win32)

# Insert two new blank rows right above the old row 47
# (TVPCannotCreateEmptyLayerImage), pushing everything else down by two.
$wsCore.Rows("47:48").Insert()

# New row 47: TVPSpecifyLayerOrBitmap
$wsCore.Range("A47").Value = "TVPSpecifyLayerOrBitmap"
$wsCore.Range("B47").Value = "Layer クラスか Bitmap クラスのオブジェクトを指定してください"
$wsCore.Range("C47").Value = "Specify Layer or Bitmap class object"

# New row 48: TVPCannotAcceptModeAuto
$wsCore.Range("A48").Value = "TVPCannotAcceptModeAuto"
$wsCore.Range("C48").Value = "Cannot accept omAuto mode"
$wsCore.Range("B48").Value = "この操作で mode に omAuto を指定することは出来ません"

# Move the view: コア共通 becomes the active sheet/tab, scrolled near the
# newly-edited rows with B49 selected; コアWin32 (previously active) is
# scrolled back up and loses the active-tab mark.
$wsWin32.Activate()
$wsWin32.Range("A13").Select()
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1

$wsCore.Activate()
$excel.ActiveWindow.ScrollRow = 28
$excel.ActiveWindow.ScrollColumn = 1
$wsCore.Range("B49").Select()
